# Update "想去人数" (number of people interested) counts for two events
# on both the "展览" sheet and the mirrored "全部类型" sheet.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 219
    $ws.Range("F4").Value = 143
}
